# Update "want to go" counts (column F) on the 展览 and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 389
$ws1.Range("F7").Value = 2444
$ws1.Range("F9").Value = 6408
$ws1.Range("F11").Value = 414
$ws1.Range("F12").Value = 26

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 389
$ws4.Range("F9").Value = 2444
$ws4.Range("F11").Value = 6408
$ws4.Range("F13").Value = 414
$ws4.Range("F15").Value = 26
